$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.910.61"

$ws.Range("D3").Value = "2.208.18"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'288.40"
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("D6").Value = "'87.33"
$ws.Range("E6").Value = "  +4.11%  "

$ws.Range("D7").Value = "'0.514"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").Value = "'30.43"
$ws.Range("E10").Value = "  +2.60%  "

$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("E12").Value = "  +2.58%  "

$ws.Range("D13").Value = "'6.43"
$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").Value = "2.555.69"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "'13.91"
$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("D16").Value = "2.204.35"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").Value = "'0.725"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "39.868.32"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").Value = "'11.59"
$ws.Range("E19").Value = "  +11.18%  "

$ws.Range("D20").Value = "0.0₃0882"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "'5.78"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "'65.38"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").Value = "'234.99"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("D26").Value = "'1.82"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'22.49"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'9.19"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'152.41"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").Value = "'31.77"
$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'4.93"
$ws.Range("E33").Value = "  +2.60%  "

$ws.Range("D34").Value = "'0.0716"
$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").Value = "'2.81"
$ws.Range("E36").Value = "  +5.54%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.0987"
$ws.Range("E38").Value = "  +2.14%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.67"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("E40").Value = "  +2.77%  "

$ws.Range("D41").Value = "'3.84"
$ws.Range("E41").Value = "  +4.47%  "

$ws.Range("D42").Value = "2.095.24"
$ws.Range("E42").Value = "  +7.94%  "

$ws.Range("D43").Value = "'2.18"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("D45").Value = "'9.93"
$ws.Range("E45").Value = "  +5.77%  "

$ws.Range("D46").Value = "'17.38"
$ws.Range("E46").Value = "  +7.99%  "

$ws.Range("D47").Value = "'2.65"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("D48").Value = "2.428.30"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("D50").Value = "'69.11"
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").Value = "'88.16"
$ws.Range("E51").Value = "  -0.67%  "
